$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new data row above the current row 27 (the "OPTIDEX -T EYE DROPS"
# row). This shifts all rows 27..46 down by one (to 28..47), carrying their
# values, styles and merged ranges with them automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Insert()

# Give the freshly inserted row the same height used for this "tall" row
# style elsewhere in the sheet.
$ws.Rows.Item(27).RowHeight = 25.5

# ---------------------------------------------------------------------------
# Fill the new row 27 with the new product line: NEXIUM 40MG 28 F.C. TAB.
# Columns C, H, L, N, O, P and Q are stored as literal text in this report
# (even though some carry a numeric-looking display format), so every value
# is written with a leading apostrophe to force text storage and preserve
# the exact formatted string (e.g. trailing zeros).
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = 21
$ws.Range("C27").Value = "'NEXIUM 40MG 28 F.C. TAB."
$ws.Range("H27").Value = "'1:0"
$ws.Range("L27").Value = "'1"
$ws.Range("N27").Value = "'488.00"
$ws.Range("P27").Value = "'3904.0000"
$ws.Range("Q27").Value = "'8:0"

# Re-create the merged cells for the new row (A:B, C:G, H:K, L:M, N:O),
# matching the pattern used by every other item row.
$ws.Range("A27:B27").Merge()
$ws.Range("C27:G27").Merge()
$ws.Range("H27:K27").Merge()
$ws.Range("L27:M27").Merge()
$ws.Range("N27:O27").Merge()

# ---------------------------------------------------------------------------
# Renumber the "#" column (A) sequentially for every item row so it keeps
# counting 1..39 after the insertion (rows 7..45 after the shift).
# ---------------------------------------------------------------------------
for ($r = 7; $r -le 45; $r++) {
    $ws.Range("A$r").Value = ($r - 6)
}

# ---------------------------------------------------------------------------
# The grand-total cell (now on row 46, previously row 45) is a plain cached
# number, not a formula - update it to include the new line's total.
# ---------------------------------------------------------------------------
$ws.Range("P46").Value = 5431.4099999999999

# ---------------------------------------------------------------------------
# The footer row (now row 47, previously row 46) keeps its text, except the
# generation timestamp which moves from 5:45 PM to 5:48 PM.
# ---------------------------------------------------------------------------
$ws.Range("A47").Value = "'Wednesday, 23 July, 2025 5:48 PM"
